$wb = $excel.ActiveWorkbook

# 1) Rename header on "Weekly Quantity" sheet
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

# 2) Rename header on "Monthly Trend" sheet
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# 3) Add new "PO Forecast" sheet after "Monthly Trend"
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45207.99999999999, 116, -69.55703518226538, 306.4058641225797),
    @(45214.99999999999, 116, -73.80345687142938, 307.9384011612991),
    @(45221.99999999999, 116, -66.9073018600035, 313.2411286945088),
    @(45228.99999999999, 115, -72.64127798695765, 300.499974978896),
    @(45235.99999999999, 115, -80.23436777929643, 304.0873095080783),
    @(45242.99999999999, 115, -74.99554244883802, 291.4892983420579),
    @(45249.99999999999, 115, -69.50092714481045, 297.2794326597208),
    @(45263.99999999999, 115, -75.7297635477822, 304.5422645775899),
    @(45305.99999999999, 115, -70.0605267014581, 320.6222266047191),
    @(45312.99999999999, 115, -84.78888788070594, 298.6454642988592),
    @(45326.99999999999, 115, -68.41849088689368, 303.0990505514989),
    @(45333.99999999999, 115, -83.45929302474377, 298.6966281807966),
    @(45347.99999999999, 115, -81.55673095482375, 306.5431865894091),
    @(45354.99999999999, 115, -79.21686808925536, 311.1666166560167),
    @(45361.99999999999, 115, -68.19878247734879, 291.2160686242877),
    @(45375.99999999999, 115, -65.5246360274392, 299.6131723095555),
    @(45382.99999999999, 115, -75.43308576145121, 301.1376967272943),
    @(45389.99999999999, 115, -83.21336732915037, 307.6670314319189),
    @(45396.99999999999, 115, -81.00940802459674, 300.0631945832149),
    @(45424.99999999999, 115, -80.61650331943859, 305.5700608050794),
    @(45431.99999999999, 115, -81.19100760485111, 298.2585154090377),
    @(45438.99999999999, 115, -83.25391638969889, 303.933106768955),
    @(45445.99999999999, 115, -73.77337692047493, 311.334273384322),
    @(45452.99999999999, 115, -75.81867133760379, 293.2183132119405),
    @(45459.99999999999, 115, -72.08598405832669, 313.058113119758),
    @(45466.99999999999, 114, -70.84714148747751, 307.5820849476381),
    @(45480.99999999999, 114, -79.30419457646708, 307.8862942567603),
    @(45487.99999999999, 114, -84.78517314654721, 286.3550354380097),
    @(45501.99999999999, 114, -66.775270581092, 297.9178315205265),
    @(45529.99999999999, 114, -84.04333631330135, 310.4636201903839),
    @(45536.99999999999, 114, -72.24668520334188, 305.2640425037719),
    @(45550.99999999999, 114, -69.56962537686417, 312.236605916861),
    @(45557.99999999999, 114, -72.96009327789751, 295.793578031733),
    @(45571.99999999999, 114, -73.99891172707852, 302.3334486171569),
    @(45634.99999999999, 114, -80.59871732591265, 293.166393057529),
    @(45641.99999999999, 114, -85.53690979605777, 296.9826938476161),
    @(45648.99999999999, 114, -70.25927156724522, 322.2142017939813),
    @(45655.99999999999, 114, -79.75154642207076, 295.0909297793536),
    @(45662.99999999999, 114, -88.66782270039646, 293.7036122723248),
    @(45669.99999999999, 114, -74.90419037072003, 294.4539642327704),
    @(45676.99999999999, 114, -65.97602952176973, 300.4939584472251),
    @(45683.99999999999, 114, -78.70686869493315, 286.0112812933309),
    @(45690.99999999999, 114, -87.42146933775921, 313.8301429405902)
)

$row = 2
foreach ($r in $data) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}

$lastRow = $row - 1

# Apply the same visual formatting used on the other two sheets:
#  - header row: bold, bordered, centered style (matches "Weekly Quantity" / "Monthly Trend" headers)
#  - column A (dates): date/time number format
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A$lastRow").PasteSpecial(-4122)

$excel.CutCopyMode = 0
